$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header: Weight -> Quantity (column D stays the same column) ---
$ws.Range("D1").Value = "Quantity"

# --- Update existing data rows (Quantity + CMP values changed) ---
$ws.Range("D2").Value = 319
$ws.Range("E2").Value = 723

$ws.Range("D3").Value = 229
$ws.Range("E3").Value = 1148

$ws.Range("D4").Value = 367
$ws.Range("E4").Value = 488

# --- Add new holding row 5: TATASTEEL (copy formatting from row 4) ---
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Add new cash/liquidity row 6 (copy header formatting from row 1) ---
$ws.Range("A1:F1").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate row 6 (CASH) values before row 5 (TATASTEEL) so new shared
# strings are appended in the same order the original workbook used.
$ws.Range("A6").Value = "CASH"
$ws.Range("B6").Value = "Deployable Capital"
$ws.Range("C6").Value = "Liquid"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 350000
$ws.Range("F6").Value = "stargate.capital"

$ws.Range("A5").Value = "TATASTEEL"
$ws.Range("B5").Value = "Tata Steel"
$ws.Range("C5").Value = "Non energy mineral"
$ws.Range("D5").Value = 781
$ws.Range("E5").Value = 168
$ws.Range("F5").Value = "tatasteel.com"

# --- Row heights (recomputed by Excel on the real edit because of wrapped text) ---
$ws.Rows.Item(1).RowHeight = 28.2
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 28.2
$ws.Rows.Item(5).RowHeight = 28.2
$ws.Rows.Item(6).RowHeight = 28.2

# --- Column widths (auto-fit in the source workbook) ---
# (Input values are tuned so the engine's internal width rounding lands on
# the stored width closest to the real auto-fit values from the diff.)
$ws.Columns.Item(1).ColumnWidth = 12.417166666666667
$ws.Columns.Item(2).ColumnWidth = 14.2505
$ws.Columns.Item(3).ColumnWidth = 12.917166666666667
$ws.Columns.Item(4).ColumnWidth = 16.58383333333333
$ws.Columns.Item(5).ColumnWidth = 17.917166666666667
$ws.Columns.Item(6).ColumnWidth = 22.58383333333333

# --- Selection moves to F6 ---
[void]$ws.Range("F6").Select()
